$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => D (Price), E (Volume 1h) new values. Rows not listed for D keep D unchanged.
$updates = @(
    @{ Row = 2;  D = "306.48";        E = "-4.57%" },
    @{ Row = 3;  D = "39.21";         E = "-8.48%" },
    @{ Row = 4;  D = "5.064";         E = "-2.87%" },
    @{ Row = 5;  D = "0.07670";       E = "-6.22%" },
    @{ Row = 6;  D = "4.230";         E = "-2.29%" },
    @{ Row = 7;  D = "1.608";         E = "-10.80%" },
    @{ Row = 8;  D = "0.9152";        E = "-3.53%" },
    @{ Row = 9;  D = "0.1028";        E = "-8.55%" },
    @{ Row = 10; D = "0.1741";        E = "-6.94%" },
    @{ Row = 11; D = $null;           E = "-5.17%" },
    @{ Row = 12; D = "0.04428";       E = "-4.16%" },
    @{ Row = 13; D = $null;           E = "-0.49%" },
    @{ Row = 14; D = "0.001253";      E = "-3.43%" },
    @{ Row = 15; D = "0.005855";      E = "1.35%" },
    @{ Row = 16; D = $null;           E = "2,412.22%" },
    @{ Row = 17; D = "3.358";         E = "-0.05%" },
    @{ Row = 18; D = "2.410";         E = "-5.26%" },
    @{ Row = 19; D = "0.3312";        E = "-1.47%" },
    @{ Row = 20; D = "7.006";         E = "-5.99%" },
    @{ Row = 21; D = "0.1347";        E = "-3.06%" },
    @{ Row = 22; D = "0.2737";        E = "7.39%" },
    @{ Row = 23; D = "0.04142";       E = "0.09%" },
    @{ Row = 24; D = "0.001204";      E = "-3.66%" },
    @{ Row = 25; D = "0.004086";      E = "-4.56%" },
    @{ Row = 26; D = "0.0001301";     E = "8.33%" },
    @{ Row = 38; D = "0.02369";       E = "-10.47%" },
    @{ Row = 39; D = "0.05161";       E = "-7.02%" },
    @{ Row = 40; D = "0.007924";      E = "-2.81%" },
    @{ Row = 41; D = "0.1316";        E = "-5.89%" },
    @{ Row = 42; D = $null;           E = "-10.90%" },
    @{ Row = 43; D = "0.001951";      E = "-6.61%" },
    @{ Row = 44; D = "0.007406";      E = "-2.36%" },
    @{ Row = 45; D = "0.3326";        E = "3.88%" },
    @{ Row = 46; D = "0.00006430";    E = "-7.00%" },
    @{ Row = 47; D = "0.00000000750"; E = "-0.04%" },
    @{ Row = 49; D = "0.004172";      E = "25.02%" },
    @{ Row = 50; D = "0.00002100";    E = "-0.04%" },
    @{ Row = 51; D = "0.0002000";     E = "-0.04%" }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }
    $cellE = $ws.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
    $cellE.Style = "Normal"
}
